$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '75.948.12'
$ws.Range('E2').Value = '  +0.49%  '
$ws.Range('D3').Value = '2.895.83'
$ws.Range('E3').Value = '  +6.45%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').Value = "'196.33"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +5.05%  '
$ws.Range('D6').Value = "'598.34"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.07%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('E8').Value = '  +2.43%  '
$ws.Range('E9').Value = '  -1.65%  '
$ws.Range('D10').Value = '2.893.35'
$ws.Range('E10').Value = '  +6.40%  '
$ws.Range('D11').Value = "'0.405"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +11.46%  '
$ws.Range('E12').Value = '  -1.26%  '
$ws.Range('E13').Value = '  +2.30%  '
$ws.Range('D14').Value = '3.427.40'
$ws.Range('E14').Value = '  +6.39%  '
$ws.Range('D15').Value = '75.890.04'
$ws.Range('E15').Value = '  +0.56%  '
$ws.Range('E16').Value = '  +0.62%  '
$ws.Range('D17').Value = "'27.37"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +1.33%  '
$ws.Range('D18').Value = '2.909.21'
$ws.Range('E18').Value = '  +7.58%  '
$ws.Range('D19').Value = "'8.91"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -4.38%  '
$ws.Range('D20').Value = "'12.60"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +3.61%  '
$ws.Range('D21').Value = "'377.55"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.11%  '
$ws.Range('D22').Value = "'2.32"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.56%  '
$ws.Range('D23').Value = "'4.16"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.11%  '
$ws.Range('D24').Value = "'71.47"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.70%  '
$ws.Range('E25').Value = '  +0.42%  '
$ws.Range('D26').Value = '3.037.32'
$ws.Range('E26').Value = '  +6.24%  '
$ws.Range('D27').Value = "'4.23"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.03%  '
$ws.Range('D28').Value = "'9.82"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +2.62%  '
$ws.Range('E29').Value = '  +10.39%  '
$ws.Range('D30').Value = "'1.00"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.03%  '
$ws.Range('E31').Value = '  -0.49%  '
$ws.Range('D32').Value = "'506.98"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -3.52%  '
$ws.Range('E33').Value = '  -1.00%  '
$ws.Range('E34').Value = '  +0.38%  '
$ws.Range('D35').Value = "'0.999"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.02%  '
$ws.Range('B36').Value = 'Monero'
$ws.Range('C36').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D36').Value = "'164.99"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +2.38%  '
$ws.Range('B37').Value = 'EthereumClassic'
$ws.Range('C37').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D37').Value = "'20.26"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +2.92%  '
$ws.Range('E38').Value = '  +1.71%  '
$ws.Range('D39').Value = "'0.114"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -4.77%  '
$ws.Range('D40').Value = "'183.49"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +5.77%  '
$ws.Range('E41').Value = '  +0.05%  '
$ws.Range('D42').Value = "'0.345"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +3.39%  '
$ws.Range('D43').Value = "'5.03"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.49%  '
$ws.Range('E44').Value = '  -2.49%  '
$ws.Range('D45').Value = "'0.0919"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +7.40%  '
$ws.Range('E46').Value = '  -0.58%  '
$ws.Range('D47').Value = "'40.37"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +3.02%  '
$ws.Range('D49').Value = "'0.581"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +5.88%  '
$ws.Range('D50').Value = "'0.669"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +11.21%  '
$ws.Range('E51').Value = '  +1.44%  '
